$d = $word.ActiveDocument

# 1. Split the "Login / Signup" list item into two separate list items:
#    "Login" (to be struck through, i.e. marked done) and "Signup" (still open).
$d.Content.Find.Execute("Login / Signup", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Login`rSignup", 2) | Out-Null

# 2. Mark completed back-end tasks as done (strikethrough).
#    After the split above, the paragraphs are, in order:
#    1 Cookbook Back-End Development
#    2 Login                              <- strike
#    3 Signup
#    4 Main Page                          <- strike
#    5 show login if no user cookie       <- strike
#    6 Show all recipes                   <- strike
#    7 Show all recipes for user          <- strike
#    8 Show all cookbooks for user
#    9 Show all recipes for a cookbook
#    10 Show recipe
#    11 Create/Edit cookbook
#    12 Create/Edit recipe
#    13 Create/Edit Ingredient
#    14 Search
$d.Paragraphs.Item(2).Range.Font.StrikeThrough = $true   # Login
$d.Paragraphs.Item(4).Range.Font.StrikeThrough = $true   # Main Page
$d.Paragraphs.Item(5).Range.Font.StrikeThrough = $true   # show login if no user cookie
$d.Paragraphs.Item(6).Range.Font.StrikeThrough = $true   # Show all recipes
$d.Paragraphs.Item(7).Range.Font.StrikeThrough = $true   # Show all recipes for user

# 3. Move the "_GoBack" bookmark so it sits at the end of the
#    "Show all recipes for a cookbook" paragraph instead of at the end of
#    the "Create/Edit Ingredient" paragraph. Re-adding a bookmark with the
#    same name ("_GoBack") automatically removes the previous one.
$pCookbook = $d.Paragraphs.Item(9)   # "Show all recipes for a cookbook"
$rCookbook = $pCookbook.Range
$rCookbookText = $d.Range($rCookbook.Start, $rCookbook.End - 1)
$d.Bookmarks.Add("_GoBack", $rCookbookText) | Out-Null
